$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 63
$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "Bhavesh"
$ws.Range("F2").Value = "bhavesh is interested in the premium package that includes advanced analytics, priority support, and additional storage capacity. She wants a detailed demo before making the decision."
$ws.Range("G2").Value = "completed"
$ws.Range("H2").Value = "'91"
$ws.Range("I2").Value = "Healthcare"
$ws.Range("J2").Value = "ABC Company"
$ws.Range("K2").Value = "Texas, USA"
$ws.Range("M2").Value = "nan`n[2025-08-25 22:45:02] The customer is interested in the premium package and wants to schedule a meeting with a representative. The meeting is scheduled for the next day at 4 p.m."
$ws.Range("N2").Value = "nan`n[2025-08-25 22:45:02] 1. Schedule a meeting with a representative, 2. Send a meeting invitation to the customer's email"
